# feat: add 2022-Q4 data
#
# 1) Insert a brand-new worksheet "2022-Q4" right after "总计" (i.e. before
#    the existing "2022-Q3" sheet), and fill it with the Q4 fund-holding data.
# 2) Add a new leading data row to the "总计" (summary) sheet describing the
#    new quarter, shifting the existing 2022-Q3 / 2022-Q2 rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the new "2022-Q4" sheet, positioned before "2022-Q3".
# ---------------------------------------------------------------------------
$q3Before = $wb.Worksheets.Item(2)
$created = $wb.Worksheets.Add($q3Before)
$created.Name = "2022-Q4"

# IMPORTANT: after Worksheets.Add() changes the sheet collection, old sheet
# references can go stale - re-resolve every sheet we need by its (now
# current) tab position.
$q4 = $wb.Worksheets.Item(2)
$q3 = $wb.Worksheets.Item(3)
$q2 = $wb.Worksheets.Item(4)

# ---- formatting: reuse the existing header / index-column style (style id
# 2 in the original file: bold, bordered, centered) by copying it over from
# the neighbouring "2022-Q3" sheet, rather than re-building it by hand.
$q3.Range("B1:H1").Copy()
$q4.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats

$q3.Range("A2").Copy()
$q4.Range("A2:A4").PasteSpecial(-4122)   # xlPasteFormats

# ---- header row
$q4.Cells.Item(1, 2).Value = "基金代码"
$q4.Cells.Item(1, 3).Value = "基金名称"
$q4.Cells.Item(1, 4).Value = "基金规模"
$q4.Cells.Item(1, 5).Value = "股票总仓位"
$q4.Cells.Item(1, 6).Value = "仓位占比"
$q4.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q4.Cells.Item(1, 8).Value = "仓位排名"

# Columns B, D, E, F, G hold numeric-looking values that must stay TEXT
# (matching the source data, e.g. fund code "006049" keeps its leading
# zero, and "4.20" keeps its trailing zero) - force a text number format
# before assigning, otherwise Excel auto-coerces them to numbers. (Union
# ranges only honour the format on their first area, so set each
# contiguous block separately.)
$q4TextRange1 = $q4.Range("B2:B4")
$q4TextRange2 = $q4.Range("D2:G4")
$q4TextRange1.NumberFormat = "@"
$q4TextRange2.NumberFormat = "@"

# ---- data rows (A = running index, H = rank are genuine numbers)
$q4.Cells.Item(2, 1).Value = 0
$q4.Cells.Item(2, 2).Value = "006049"
$q4.Cells.Item(2, 3).Value = "恒越研究精选混合A/B"
$q4.Cells.Item(2, 4).Value = "4.23"
$q4.Cells.Item(2, 5).Value = "88.62"
$q4.Cells.Item(2, 6).Value = "4.20"
$q4.Cells.Item(2, 7).Value = "0.1777"
$q4.Cells.Item(2, 8).Value = 4

$q4.Cells.Item(3, 1).Value = 1
$q4.Cells.Item(3, 2).Value = "012846"
$q4.Cells.Item(3, 3).Value = "恒越蓝筹精选混合"
$q4.Cells.Item(3, 4).Value = "5.90"
$q4.Cells.Item(3, 5).Value = "86.65"
$q4.Cells.Item(3, 6).Value = "2.76"
$q4.Cells.Item(3, 7).Value = "0.1628"
$q4.Cells.Item(3, 8).Value = 7

$q4.Cells.Item(4, 1).Value = 2
$q4.Cells.Item(4, 2).Value = "007192"
$q4.Cells.Item(4, 3).Value = "恒越研究精选混合C"
$q4.Cells.Item(4, 4).Value = "3.19"
$q4.Cells.Item(4, 5).Value = "88.62"
$q4.Cells.Item(4, 6).Value = "4.20"
$q4.Cells.Item(4, 7).Value = "0.1340"
$q4.Cells.Item(4, 8).Value = 4

# Restore the originally-selected tab (adding a sheet above made the new
# sheet active; the source file had "2022-Q2" selected).
$q2.Activate()

# ---------------------------------------------------------------------------
# Step 2: insert the new Q4 summary row into "总计", pushing the existing
# 2022-Q3 / 2022-Q2 rows down one row each.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(1)

# Preserve the bordered/bold "index column" style (style id 2) on the row
# that is about to become row 4 as well, by copying format from A2 first.
$total.Range("A2").Copy()
$total.Range("A3:A4").PasteSpecial(-4122)   # xlPasteFormats

# Row 4 <- old row 3 (2022-Q2)
$total.Cells.Item(4, 2).Value = "2022-Q2"
$total.Cells.Item(4, 3).Value = 1
$total.Cells.Item(4, 4).Value = 0
$total.Cells.Item(4, 1).Value = 2

# Row 3 <- old row 2 (2022-Q3)
$total.Cells.Item(3, 2).Value = "2022-Q3"
$total.Cells.Item(3, 3).Value = 4
$total.Cells.Item(3, 4).Value = 0.08
$total.Cells.Item(3, 1).Value = 1

# Row 2 <- new 2022-Q4 figures
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.47
$total.Cells.Item(2, 1).Value = 0
